$wb = $excel.ActiveWorkbook

# Sheet "1" (index 2): add "Score" header in C1, update selection to C1
$ws = $wb.Worksheets.Item(2)
$ws.Activate()
$ws.Range("C1").Value = "Score"
$ws.Range("C1").Select()

# Sheet "2" (index 3): add "Score" header in C1, update selection to C1
$ws = $wb.Worksheets.Item(3)
$ws.Activate()
$ws.Range("C1").Value = "Score"
$ws.Range("C1").Select()

# Sheet "3" (index 4): add "Score" header in C1, update selection to C1
$ws = $wb.Worksheets.Item(4)
$ws.Activate()
$ws.Range("C1").Value = "Score"
$ws.Range("C1").Select()

# Sheet "4" (index 5): add "Score" header in C1, update selection to C1
$ws = $wb.Worksheets.Item(5)
$ws.Activate()
$ws.Range("C1").Value = "Score"
$ws.Range("C1").Select()

# Sheet "5" (index 6): add "Score" header in D1, update selection to D1
$ws = $wb.Worksheets.Item(6)
$ws.Activate()
$ws.Range("D1").Value = "Score"
$ws.Range("D1").Select()

# Sheet "6" (index 7): add "Score" header in C1, update selection to C1
$ws = $wb.Worksheets.Item(7)
$ws.Activate()
$ws.Range("C1").Value = "Score"
$ws.Range("C1").Select()

# Sheet "7" (index 8): add "Score" header in C1, update selection to H11, make this the active tab
$ws = $wb.Worksheets.Item(8)
$ws.Activate()
$ws.Range("C1").Value = "Score"
$ws.Range("H11").Select()
